$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRED Graph")

# Update existing GDP values for rows 12-30 (years 2001-2019)
$ws.Range("B12").Value = 159243.61199999999
$ws.Range("B13").Value = 162437.21100000001
$ws.Range("B14").Value = 167612.44699999999
$ws.Range("B15").Value = 175236.78
$ws.Range("B16").Value = 192044.47399999999
$ws.Range("B17").Value = 206996.83300000001
$ws.Range("B18").Value = 227215.084
$ws.Range("B19").Value = 235857.435
$ws.Range("B20").Value = 231980.182
$ws.Range("B21").Value = 240147.40900000001
$ws.Range("B22").Value = 251468.13200000001
$ws.Range("B23").Value = 269458.587
$ws.Range("B24").Value = 285068.19400000002
$ws.Range("B25").Value = 302314.75599999999
$ws.Range("B26").Value = 320665.74800000002
$ws.Range("B27").Value = 336256.81599999999
$ws.Range("B28").Value = 360610.815
$ws.Range("B29").Value = 393634.70199999999
$ws.Range("B30").Value = 419475.451

# Add new row 31 for year 2020
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 426939.96299999999
$ws.Range("B31").NumberFormat = "0.000"

# Update selection to reflect the edit state: the user's cursor had moved to
# the new first-empty row (A32) after adding the 2020 observation, with the
# full data columns (A:B) selected.
$null = $ws.Range("A1:B1048576").Select()
